$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (pt_max) values from 50 to 70 for rows 2 through 12
$ws.Range("E2:E12").Value = 70

# Update the active cell selection to E19
$ws.Range("E19").Select()
